$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the Jar info rows (14-16), mirroring the "Total" row's boxed style.
$ws.Range("A14").Value = "Jar Description"
$ws.Range("C14").Value = "Diameter"

$ws.Range("A15").Value = "Jar Inner Opening "
$ws.Range("C15").Value = "2.04`" (52 mm)"

$ws.Range("A16").Value = "Jar Outer Opening "
$ws.Range("C16").Value = "2.32`" (59 mm)"

# Borders to match the existing "Total" row box styling.
$ws.Range("A14").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("A14:C16").Borders.Item(8).LineStyle = 1 # xlEdgeTop on whole range (approx)

$ws.Range("E11").Select()
